$d = $word.ActiveDocument

# ---------- helpers ----------------------------------------------------
function Set-ParaText($idx1based, $newtext) {
    $p = $d.Paragraphs($idx1based)
    $full = $p.Range
    $r = $d.Range($full.Start, $full.End - 1)
    $r.Text = $newtext
}

function Split-RunAt($pos) {
    # Forces a run boundary at an absolute character offset by wrapping a
    # throw-away bookmark around the (zero-length) point and immediately
    # deleting it again -- the run split persists after the bookmark goes.
    $b = $d.Bookmarks.Add("TEMPSPLIT", $d.Range($pos, $pos))
    $d.Bookmarks("TEMPSPLIT").Delete()
}

# ---------- 1. drop the _GoBack bookmark from the title paragraph ------
$d.Bookmarks("_GoBack").Delete()

# ---------- 2. "Register " paragraph gains extra runs ------------------
$pReg = $d.Paragraphs(4)
$regFull = $pReg.Range
$regParaStart = $regFull.Start
$regInsPos = $regFull.End - 1
$regIp = $d.Range($regInsPos, $regInsPos)
$regIp.InsertAfter("– url + request method /users/re  ")
Split-RunAt ($regParaStart + 9)
Split-RunAt ($regParaStart + 11)
Split-RunAt ($regParaStart + 14)
Split-RunAt ($regParaStart + 41)

# ---------- 3. "Watch video" / "Subscribe to  user " swap ---------------
Set-ParaText 7 "Subscribe to  user "
Set-ParaText 8 "Unsubscribe to user"

# ---------- 4. re-home the _GoBack bookmark on the spacer paragraph ----
$pSpacer = $d.Paragraphs(10)
$spacerFull = $pSpacer.Range
$spacerEndPos = $spacerFull.End - 1
$zPoint = $d.Range($spacerEndPos, $spacerEndPos)
$zPoint.InsertAfter("Z")
$zRange = $d.Range($spacerEndPos, $spacerEndPos + 1)
$d.Bookmarks.Add("_GoBack", $zRange) | Out-Null
$d.Range($spacerEndPos, $spacerEndPos + 1).Text = ""

# ---------- 5. simple text swaps ---------------------------------------
Set-ParaText 12 "Add video"
Set-ParaText 15 "Get by id"

# ---------- 6. "Find video by title" -> two runs -----------------------
$pFind = $d.Paragraphs(16)
$findFull = $pFind.Range
$findStart = $findFull.Start
$findR = $d.Range($findFull.Start, $findFull.End - 1)
$findR.Text = "Get all by title"
Split-RunAt ($findStart + 10)

# ---------- 7. "Create " -> "add" ---------------------------------------
Set-ParaText 20 "add"

# ---------- 8. new "Remove video from playlist" bullet -----------------
$pAddVideo = $d.Paragraphs(21)
$pAddVideo.Range.InsertParagraphAfter()
$pNewPlaylist = $d.Paragraphs(22)
$pNewPlaylist.Range.Text = "Remove video from playlist"

# ---------- 9. new "Edit" / "Delete" bullets + trailing blank para -----
$pDislike2 = $d.Paragraphs(30)
$pDislike2.Range.InsertParagraphAfter()
$pEdit = $d.Paragraphs(31)
$pEdit.Range.Text = "Edit"

$pEdit2 = $d.Paragraphs(31)
$pEdit2.Range.InsertParagraphAfter()
$pDelete = $d.Paragraphs(32)
$pDelete.Range.Text = "Delete"

$pDelete2 = $d.Paragraphs(32)
$pDelete2.Range.InsertParagraphAfter()
$pBlankNew = $d.Paragraphs(33)
$pBlankNew.Range.ListFormat.RemoveNumbers()

Write-Output "edit complete"
